# Apply odds updates to the FlashScore weekly games sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (Los Andes - Alvarado)
$ws.Range("N3").Value = 2.88
$ws.Range("O3").Value = 1.4

# Row 4 (CA Estudiantes - Temperley)
$ws.Range("O4").Value = 1.33
$ws.Range("R4").Value = 2.63
$ws.Range("S4").Value = 1.44

# Row 9 (Independiente FBC - River Plate)
$ws.Range("G9").Value = 1.85
$ws.Range("H9").Value = 3.4
$ws.Range("I9").Value = 3.95
$ws.Range("L9").Value = 1.32
$ws.Range("M9").Value = 2.85
$ws.Range("N9").Value = 1.93
$ws.Range("O9").Value = 1.7
$ws.Range("P9").Value = 1.42
$ws.Range("Q9").Value = 2.45
$ws.Range("R9").Value = 1.82
$ws.Range("S9").Value = 1.78
$ws.Range("T9").Value = 6.5
$ws.Range("U9").Value = 8.25
$ws.Range("V9").Value = 8.5
$ws.Range("X9").Value = 15.5
$ws.Range("Y9").Value = 30
$ws.Range("Z9").Value = 9
$ws.Range("AA9").Value = 6.6
$ws.Range("AB9").Value = 16
$ws.Range("AC9").Value = 80
$ws.Range("AD9").Value = 10.75
$ws.Range("AE9").Value = 21
$ws.Range("AG9").Value = 60
$ws.Range("AH9").Value = 37
$ws.Range("AJ9").Value = 700

# Row 11 (Cheongju - Asan)
$ws.Range("L11").Value = 1.33
$ws.Range("M11").Value = 3.25
$ws.Range("N11").Value = 2.05
$ws.Range("O11").Value = 1.75
$ws.Range("W11").Value = 41
$ws.Range("X11").Value = 29
$ws.Range("Y11").Value = 41
$ws.Range("AD11").Value = 7
$ws.Range("AE11").Value = 9.5
$ws.Range("AH11").Value = 17
$ws.Range("AJ11").Value = 301
